# "change plate for RFLP analysis"
# The 008 plate has been run/processed, so its tab is renamed with a
# "DONE " prefix and it becomes the active sheet (the next plate to work
# on, 009, keeps the cursor position that was left on it, and the sheet
# that used to be selected - 013 - goes back to being a plain tab).

$wb = $excel.ActiveWorkbook

# Rename "PocHistone RLFP 008" -> "DONE PocHistone RLFP 008"
$doneSheet = $wb.Worksheets.Item("PocHistone RLFP 008")
$doneSheet.Name = "DONE PocHistone RLFP 008"

# Update the leftover selection on that sheet and make it the active tab
[void]$doneSheet.Select()
[void]$doneSheet.Range("E14").Select()

# That sheet is also switched to landscape / narrower margins for printing
$doneSheet.PageSetup.Orientation = 2   # xlLandscape
$doneSheet.PageSetup.PaperSize = 9     # xlPaperA4
$doneSheet.PageSetup.LeftMargin = $excel.InchesToPoints(0.25)
$doneSheet.PageSetup.RightMargin = $excel.InchesToPoints(0.25)

# Move the selection on "PocHistone RLFP 009" one cell to the right
$nextSheet = $wb.Worksheets.Item("PocHistone RLFP 009")
[void]$nextSheet.Range("C10").Select()

# Every other sheet just picks up the default page setup (A4 portrait)
# that was applied across the workbook.
foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -ne "DONE PocHistone RLFP 008") {
        $ws.PageSetup.Orientation = 1  # xlPortrait
        $ws.PageSetup.PaperSize = 9    # xlPaperA4
    }
}

[void]$doneSheet.Activate()
